# Daily attendance processing - 2026-01-04 19:07:16
# Swap the order of "System" and the user email in column G ("Recorded By")
# from "System, <email>" to "<email>, System" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)
        $cell.Value = "$rest, System"
    }
}
